# Add "NA" values under the duplicate_image_filename column (column E)
# for the data rows 2 through 21.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 21; $row++) {
    $ws.Range("E$row").Value = "NA"
}
